$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new order data
$ws.Range("A2").Value = 674
$ws.Range("B2").Value = 45763.22928240741
$ws.Range("C2").Value = "Ajay Francis Anchan"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 580
$ws.Range("F2").Value = 14.25
$ws.Range("G2").Value = 14.25
$ws.Range("H2").Value = 28.5
$ws.Range("I2").Value = "Oreo Shake (x1), Butterscotch Lassi (x4), Chicken Wrap (x4), Banana Shake (x1)"

# Remove the now-unneeded third row entirely
$ws.Rows("3:3").Delete()
